$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value2 = 'ACT-GEN-CAP-011'
$ws.Range('B2').Value2 = 'GEN-CAP-011'
$ws.Range('C2').Value2 = '[Generated] M28 Governance Orchestration'
$ws.Range('F2').Value2 = 'CAP-CAP-011'
$ws.Range('J2').Value2 = 'BOUNDED'
$ws.Range('N2').Value2 = 'CAP-011'
$ws.Range('A3').Value2 = 'ACT-GEN-CAP-016'
$ws.Range('B3').Value2 = 'GEN-CAP-016'
$ws.Range('C3').Value2 = '[Generated] Skill System (M2/M3)'
$ws.Range('F3').Value2 = 'CAP-CAP-016'
$ws.Range('N3').Value2 = 'CAP-016'
$ws.Range('A4').Value2 = 'ACT-GEN-CAP-002'
$ws.Range('B4').Value2 = 'GEN-CAP-002'
$ws.Range('C4').Value2 = '[Generated] Cost Simulation V2'
$ws.Range('F4').Value2 = 'CAP-CAP-002'
$ws.Range('N4').Value2 = 'CAP-002'
$ws.Range('A5').Value2 = 'ACT-GEN-CAP-008'
$ws.Range('B5').Value2 = 'GEN-CAP-008'
$ws.Range('C5').Value2 = '[Generated] M12 Multi-Agent Orchestration'
$ws.Range('F5').Value2 = 'CAP-CAP-008'
$ws.Range('K5').Value2 = 'READ'
$ws.Range('N5').Value2 = 'CAP-008'
$ws.Range('A6').Value2 = 'ACT-GEN-CAP-012'
$ws.Range('B6').Value2 = 'GEN-CAP-012'
$ws.Range('C6').Value2 = '[Generated] M4 Workflow Engine'
$ws.Range('F6').Value2 = 'CAP-CAP-012'
$ws.Range('J6').Value2 = 'STRICT'
$ws.Range('N6').Value2 = 'CAP-012'
$ws.Range('A7').Value2 = 'ACT-GEN-CAP-020'
$ws.Range('B7').Value2 = 'GEN-CAP-020'
$ws.Range('C7').Value2 = '[Generated] CLI Execution'
$ws.Range('F7').Value2 = 'CAP-CAP-020'
$ws.Range('J7').Value2 = 'BOUNDED'
$ws.Range('N7').Value2 = 'CAP-020'
$ws.Range('A8').Value2 = 'ACT-GEN-CAP-021'
$ws.Range('B8').Value2 = 'GEN-CAP-021'
$ws.Range('C8').Value2 = '[Generated] SDK Execution'
$ws.Range('F8').Value2 = 'CAP-CAP-021'
$ws.Range('N8').Value2 = 'CAP-021'
$ws.Range('A9').Value2 = 'ACT-GEN-CAP-001'
$ws.Range('B9').Value2 = 'GEN-CAP-001'
$ws.Range('C9').Value2 = '[Generated] Execution Replay & Activity'
$ws.Range('F9').Value2 = 'CAP-CAP-001'
$ws.Range('H9').Value2 = 'EVIDENCE'
$ws.Range('I9').Value2 = 'OBSERVE'
$ws.Range('J9').Value2 = 'STRICT'
$ws.Range('N9').Value2 = 'CAP-001'
$ws.Range('A16').Value2 = 'INC-GEN-CAP-002'
$ws.Range('B16').Value2 = 'GEN-CAP-002'
$ws.Range('C16').Value2 = '[Generated] Cost Simulation V2'
$ws.Range('F16').Value2 = 'CAP-CAP-002'
$ws.Range('J16').Value2 = 'BOUNDED'
$ws.Range('N16').Value2 = 'CAP-002'
$ws.Range('A17').Value2 = 'INC-GEN-CAP-005'
$ws.Range('B17').Value2 = 'GEN-CAP-005'
$ws.Range('C17').Value2 = '[Generated] Founder Console'
$ws.Range('F17').Value2 = 'CAP-CAP-005'
$ws.Range('I17').Value2 = 'CONTROL'
$ws.Range('K17').Value2 = 'GOVERN'
$ws.Range('N17').Value2 = 'CAP-005'
$ws.Range('A18').Value2 = 'INC-GEN-CAP-009'
$ws.Range('B18').Value2 = 'GEN-CAP-009'
$ws.Range('C18').Value2 = '[Generated] M19 Policy Engine'
$ws.Range('F18').Value2 = 'CAP-CAP-009'
$ws.Range('I18').Value2 = 'ACT'
$ws.Range('K18').Value2 = 'WRITE'
$ws.Range('N18').Value2 = 'CAP-009'
$ws.Range('A19').Value2 = 'INC-GEN-CAP-021'
$ws.Range('B19').Value2 = 'GEN-CAP-021'
$ws.Range('C19').Value2 = '[Generated] SDK Execution'
$ws.Range('F19').Value2 = 'CAP-CAP-021'
$ws.Range('K19').Value2 = 'READ'
$ws.Range('N19').Value2 = 'CAP-021'
$ws.Range('A20').Value2 = 'INC-GEN-CAP-001'
$ws.Range('B20').Value2 = 'GEN-CAP-001'
$ws.Range('C20').Value2 = '[Generated] Execution Replay & Activity'
$ws.Range('F20').Value2 = 'CAP-CAP-001'
$ws.Range('H20').Value2 = 'EVIDENCE'
$ws.Range('I20').Value2 = 'OBSERVE'
$ws.Range('J20').Value2 = 'STRICT'
$ws.Range('N20').Value2 = 'CAP-001'
$ws.Range('A21').Value2 = 'LOG-GEN-CAP-002'
$ws.Range('B21').Value2 = 'GEN-CAP-002'
$ws.Range('C21').Value2 = '[Generated] Cost Simulation V2'
$ws.Range('F21').Value2 = 'CAP-CAP-002'
$ws.Range('J21').Value2 = 'BOUNDED'
$ws.Range('N21').Value2 = 'CAP-002'
$ws.Range('A22').Value2 = 'LOG-GEN-CAP-021'
$ws.Range('B22').Value2 = 'GEN-CAP-021'
$ws.Range('C22').Value2 = '[Generated] SDK Execution'
$ws.Range('F22').Value2 = 'CAP-CAP-021'
$ws.Range('K22').Value2 = 'READ'
$ws.Range('N22').Value2 = 'CAP-021'
$ws.Range('A23').Value2 = 'LOG-GEN-CAP-001'
$ws.Range('B23').Value2 = 'GEN-CAP-001'
$ws.Range('C23').Value2 = '[Generated] Execution Replay & Activity'
$ws.Range('F23').Value2 = 'CAP-CAP-001'
$ws.Range('H23').Value2 = 'EVIDENCE'
$ws.Range('I23').Value2 = 'OBSERVE'
$ws.Range('J23').Value2 = 'STRICT'
$ws.Range('N23').Value2 = 'CAP-001'
$ws.Range('A28').Value2 = 'POL-GEN-CAP-005'
$ws.Range('B28').Value2 = 'GEN-CAP-005'
$ws.Range('C28').Value2 = '[Generated] Founder Console'
$ws.Range('F28').Value2 = 'CAP-CAP-005'
$ws.Range('I28').Value2 = 'CONTROL'
$ws.Range('J28').Value2 = 'BOUNDED'
$ws.Range('K28').Value2 = 'GOVERN'
$ws.Range('N28').Value2 = 'CAP-005'
$ws.Range('A29').Value2 = 'POL-GEN-CAP-011'
$ws.Range('B29').Value2 = 'GEN-CAP-011'
$ws.Range('C29').Value2 = '[Generated] M28 Governance Orchestration'
$ws.Range('F29').Value2 = 'CAP-CAP-011'
$ws.Range('I29').Value2 = 'ACT'
$ws.Range('K29').Value2 = 'WRITE'
$ws.Range('N29').Value2 = 'CAP-011'
$ws.Range('A30').Value2 = 'POL-GEN-CAP-003'
$ws.Range('B30').Value2 = 'GEN-CAP-003'
$ws.Range('C30').Value2 = '[Generated] Policy Proposals'
$ws.Range('F30').Value2 = 'CAP-CAP-003'
$ws.Range('H30').Value2 = 'SUBSTRATE'
$ws.Range('I30').Value2 = 'EXPLAIN'
$ws.Range('K30').Value2 = 'READ'
$ws.Range('N30').Value2 = 'CAP-003'
$ws.Range('A31').Value2 = 'POL-GEN-CAP-016'
$ws.Range('B31').Value2 = 'GEN-CAP-016'
$ws.Range('C31').Value2 = '[Generated] Skill System (M2/M3)'
$ws.Range('F31').Value2 = 'CAP-CAP-016'
$ws.Range('H31').Value2 = 'ACTION'
$ws.Range('I31').Value2 = 'ACT'
$ws.Range('K31').Value2 = 'WRITE'
$ws.Range('N31').Value2 = 'CAP-016'
$ws.Range('A32').Value2 = 'POL-GEN-CAP-001'
$ws.Range('B32').Value2 = 'GEN-CAP-001'
$ws.Range('C32').Value2 = '[Generated] Execution Replay & Activity'
$ws.Range('F32').Value2 = 'CAP-CAP-001'
$ws.Range('H32').Value2 = 'EVIDENCE'
$ws.Range('I32').Value2 = 'OBSERVE'
$ws.Range('J32').Value2 = 'STRICT'
$ws.Range('K32').Value2 = 'READ'
$ws.Range('N32').Value2 = 'CAP-001'
